$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "'29.571.89"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "'1.913.60"
$ws.Range("E3").Value = "  -0.09%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "'325.65"
$ws.Range("E5").Value = "  -0.49%  "

# Row 6: 'USDC' -> 'USDC'
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.57%  "

# Row 7: 'XRP' -> 'XRP'
$ws.Range("D7").Value = "'0.4832"
$ws.Range("E7").Value = "  +0.44%  "

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").Value = "'0.4079"
$ws.Range("E8").Value = "  -0.52%  "

# Row 9: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D9").Value = "'0.08174"

# Row 10: 'Polygon' -> 'Polygon'
$ws.Range("D10").Value = "'1.012"
$ws.Range("E10").Value = "  +0.00%  "

# Row 11: 'Solana' -> 'Solana'
$ws.Range("E11").Value = "  +4.87%  "

# Row 12: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D12").Value = "'1.919.63"
$ws.Range("E12").Value = "  +1.03%  "

# Row 13: 'Polkadot' -> 'Polkadot'
$ws.Range("D13").Value = "'6.029"
$ws.Range("E13").Value = "  +1.27%  "

# Row 14: 'Chainlink' -> 'Chainlink'
$ws.Range("D14").Value = "'7.118"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15: 'Litecoin' -> 'Litecoin'
$ws.Range("D15").Value = "'90.48"
$ws.Range("E15").Value = "  +0.92%  "

# Row 16: 'TRON' -> 'TRON'
$ws.Range("D16").Value = "'0.06793"
$ws.Range("E16").Value = "  +2.81%  "

# Row 17: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("E17").Value = "  +0.68%  "

# Row 18: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D18").Value = "'0.00001042"
$ws.Range("E18").Value = "  +0.91%  "

# Row 19: 'Avalanche' -> 'Avalanche'
$ws.Range("D19").Value = "'17.72"
$ws.Range("E19").Value = "  -0.18%  "

# Row 20: 'Dai' -> 'Dai'
$ws.Range("E20").Value = "  +0.54%  "

# Row 21: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D21").Value = "'29.576.93"
$ws.Range("E21").Value = "  +0.18%  "

# Row 22: 'Uniswap' -> 'Uniswap'
$ws.Range("D22").Value = "'5.621"
$ws.Range("E22").Value = "  +1.20%  "

# Row 23: 'Cosmos' -> 'Cosmos'
$ws.Range("D23").Value = "'11.81"
$ws.Range("E23").Value = "  +2.30%  "

# Row 24: 'Toncoin' -> 'Toncoin'
$ws.Range("D24").Value = "'2.173"
$ws.Range("E24").Value = "  -1.38%  "

# Row 25: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D25").Value = "'2.155.65"
$ws.Range("E25").Value = "  +2.03%  "

# Row 26: 'Monero' -> 'Monero'
$ws.Range("D26").Value = "'154.66"
$ws.Range("E26").Value = "  +0.54%  "

# Row 27: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D27").Value = "'20.06"
$ws.Range("E27").Value = "  +0.97%  "

# Row 28: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D28").Value = "'6.334"
$ws.Range("E28").Value = "  +9.43%  "

# Row 29: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D29").Value = "'2.106"
$ws.Range("E29").Value = "  -1.42%  "

# Row 30: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D30").Value = "'119.70"
$ws.Range("E30").Value = "  +1.58%  "

# Row 31: 'ImmutableX' -> 'ImmutableX'
$ws.Range("E31").Value = "  -2.59%  "

# Row 32: 'Stellar' -> 'Stellar'
$ws.Range("D32").Value = "'0.09580"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33: 'Filecoin' -> 'Filecoin'
$ws.Range("D33").Value = "'5.554"
$ws.Range("E33").Value = "  +2.96%  "

# Row 34: 'ARBITRUM' -> 'HuobiToken'
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.558"
$ws.Range("E34").Value = "  -0.40%  "

# Row 35: 'HuobiToken' -> 'ARBITRUM'
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.396"
$ws.Range("E35").Value = "  -2.01%  "

# Row 36: 'VeChain' -> 'VeChain'
$ws.Range("D36").Value = "'0.02270"
$ws.Range("E36").Value = "  +0.70%  "

# Row 37: 'Hedera' -> 'Hedera'
$ws.Range("E37").Value = "  +0.21%  "

# Row 38: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("D38").Value = "'1.174"
$ws.Range("E38").Value = "  +0.00%  "

# Row 39: 'TheSandbox' -> 'Aptos'
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'10.82"
$ws.Range("E39").Value = "  +6.55%  "

# Row 40: 'Aptos' -> 'TheSandbox'
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5947"
$ws.Range("E40").Value = "  +1.02%  "

# Row 41: 'FraxShare' -> 'FraxShare'
$ws.Range("D41").Value = "'7.941"
$ws.Range("E41").Value = "  -4.84%  "

# Row 42: 'Algorand' -> 'Algorand'
$ws.Range("D42").Value = "'0.1856"
$ws.Range("E42").Value = "  +0.57%  "

# Row 43: 'RenderToken' -> 'RenderToken'
$ws.Range("D43").Value = "'2.464"
$ws.Range("E43").Value = "  +0.06%  "

# Row 44: 'WEMIXToken' -> 'WEMIXToken'
$ws.Range("D44").Value = "'1.283"
$ws.Range("E44").Value = "  -0.76%  "

# Row 45: 'Cronos' -> 'Cronos'
$ws.Range("E45").Value = "  -3.97%  "

# Row 46: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D46").Value = "'12.39"
$ws.Range("E46").Value = "  +2.11%  "

# Row 47: 'Decentraland' -> 'Decentraland'
$ws.Range("D47").Value = "'0.5577"
$ws.Range("E47").Value = "  +0.50%  "

# Row 48: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D48").Value = "'1.954"
$ws.Range("E48").Value = "  +1.14%  "

# Row 49: 'Quant' -> 'Quant'
$ws.Range("D49").Value = "'115.25"
$ws.Range("E49").Value = "  +1.44%  "

# Row 50: 'Aave' -> 'Aave'
$ws.Range("D50").Value = "'72.79"
$ws.Range("E50").Value = "  +1.68%  "

# Row 51: 'EOS' -> 'EOS'
$ws.Range("E51").Value = "  +1.96%  "

